$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "RUNMANAGER" (sheet1): add a new test-scenario row (giveFeedbackOrderMcPizza)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("RUNMANAGER")

# Clone formatting of the row above into the new row, then clone the text-as-
# string "1" values for count/priority (D/E) so they stay shared-string "1"
# rather than becoming numeric.
$ws1.Range("A3:E3").Copy()
$ws1.Range("A4:E4").PasteSpecial(-4122)
$ws1.Range("D3:E3").Copy()
$ws1.Range("D4:E4").PasteSpecial(-4163)
$ws1.Rows.Item(4).RowHeight = 28.5

$ws1.Range("A4").Value = "giveFeedbackOrderMcPizza"
$ws1.Range("B4").Value = "Submit Feedback to McPIizza bot"
$ws1.Range("C4").Value = "yes"

$ws1.Range("B4").Select()

# ---------------------------------------------------------------------------
# Sheet "DATA" (sheet2): split the feedback scenario out of rows 2/3 into its
# own row 5, as a dedicated "Feedback Test" run.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("DATA")

# New hyperlink first (picks up rId4 / default style), then the formats from
# row 3 get pasted on top so the cell keeps the plain "s=5" look of F2/F3/F4.
$ws2.Hyperlinks.Add($ws2.Range("F5"), "mailto:amitnitishkumar@gmail.com")

$ws2.Range("A3:O3").Copy()
$ws2.Range("A5:O5").PasteSpecial(-4122)
$ws2.Range("C3:O3").Copy()
$ws2.Range("C5:O5").PasteSpecial(-4163)
$ws2.Rows.Item(5).RowHeight = 28.5

$ws2.Range("A5").Value = "giveFeedbackOrderMcPizza"
$ws2.Range("B5").Value = "Feedback Test"

# The feedback-pick value no longer belongs on the combined rows 2/3 - it now
# only lives on the new dedicated feedback row (5).
$ws2.Range("N2").ClearContents()
$ws2.Range("N3").ClearContents()

$ws2.Range("A2").Select()
